$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = -7.280199999999993
    "C7"  = -13.961
    "B10" = 5.1815
    "B12" = 4.964600000000002
    "C15" = -14.42839999999999
    "B18" = 6.678699999999995
    "D18" = -9.035899999999993
    "D19" = -9.05979999999999
    "C20" = -12.1368
    "D27" = -9.0053
    "C29" = -11.52640000000001
    "C30" = -12.688
    "C31" = -13.0431
    "B37" = 8.811100000000001
    "C40" = -13.2629
    "D42" = -8.944299999999997
    "D44" = -7.589300000000001
    "D47" = -7.575000000000002
    "B55" = 6.586899999999995
    "D58" = -8.17709999999999
    "B68" = 5.945399999999998
    "C68" = -11.74250000000001
    "D73" = -7.895899999999999
    "C76" = -12.7209
    "B77" = 9.332000000000006
    "B78" = 9.608800000000004
    "C87" = -14.49149999999999
    "C88" = -13.76869999999999
    "D95" = -7.7093
    "C96" = -12.81350000000001
    "C98" = -12.2257
    "C101" = -13.383
    "D101" = -8.132999999999997
    "C102" = -13.14940000000001
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
